# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps. These columns reuse the same shared-string entries across
# rows whose recorded timestamps happen to coincide (rows 2 and 5 on each
# language sheet), so every cell pointing at that value moves together.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-17 14:16:21"
$zhcn.Range("H2").Value = "2016-03-17 14:16:38"
$zhcn.Range("E5").Value = "2016-03-17 14:16:21"
$zhcn.Range("H5").Value = "2016-03-17 14:16:38"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-17 14:16:24"
$dede.Range("H2").Value = "2016-03-17 14:16:45"
$dede.Range("E5").Value = "2016-03-17 14:16:24"
$dede.Range("H5").Value = "2016-03-17 14:16:45"
